$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete trailing rows 24 and 25 (content moved up into 22-23),
# shrinking the sheet from A1:C25 to A1:C23.
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# Rewrite rows 10, 12-23 with the restructured content/layout
# Row 10
$ws.Cells.Item(10,1).Value = "Objetivos:"
$ws.Cells.Item(10,2).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Cells.Item(10,3).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item(11,1).Value = "Objectives:"
$ws.Cells.Item(11,2).Value = "The course aims to cover the concepts of biochemistry, cell structural organization and molecular composition; to understand the importance of organic compounds and cellular metabolism; and to enable students to acquire a specialised knowledge for further disciplines of the undergraduate program in Chemical Engineering and Industrial Chemical."
$ws.Cells.Item(11,3).Value = "The course aims to cover the concepts of biochemistry, cell structural organization and molecular composition; to understand the importance of organic compounds and cellular metabolism; and to enable students to acquire a specialised knowledge for further disciplines of the undergraduate program in Chemical Engineering and Industrial Chemical."
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item(12,1).Value = "Docentes responsáveis:"
$ws.Cells.Item(12,2).ClearContents()
$ws.Cells.Item(12,3).ClearContents()

# Row 13
$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "01/01/2019"
$ws.Cells.Item(13,3).Value = "01/01/2019"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Cells.Item(14,2).Value = "Physical and chemical properties of biomolecules and their levels of organization. Overview of DNA-based information technologies and protein metabolism, overview of glucose metabolism, anaerobic metabolism, oxidative metabolism of electron-transfer reaction, oxidative phosphorylation, photosynthesis."
$ws.Cells.Item(14,3).Value = "Physical and chemical properties of biomolecules and their levels of organization. Overview of DNA-based information technologies and protein metabolism, overview of glucose metabolism, anaerobic metabolism, oxidative metabolism of electron-transfer reaction, oxidative phosphorylation, photosynthesis."
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Cells.Item(15,1).Value = "Programa:"
$ws.Cells.Item(15,2).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Cells.Item(15,3).Value = "8711290 - Elisson Antônio da Costa Romanel"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Cells.Item(16,2).Value = "Basic biochemistry: Solvent properties of water, acids and bases, titration curves, buffer solutions. Aminoacids: three-dimensional structure, structure-property relationship, isoelectric point, electrophoresis, peptide bond, small peptides with physiological activity. Proteins: structure levels, structural irregularities, denaturation and renaturation, isolation, characterization and quantification. Enzymes: general concepts and mechanisms of action, cofactors and coenzymes, kinetics of michaelian enzymes, equilibrium and reaction rate, kinetic parameters and their applications, allosteric enzymes, regulation and inhibition of enzymes, general mechanisms of enzymatic reactions. Carbohydrates: classification of monosaccharides, cyclic structure and isomerism, chemical properties, disaccharides, homopolysaccharides and their structural and energy storage functions, heteropolysaccharides, glycoproteins and glycolipids. Lipids: fatty acids, triacylglycerides, phospholipids, sphingolipids and cholesterol. Biological membranes: fluid mosaic model, transport, selective permeability, passive process and active transport of biomolecules and/or ions. Molecular Biochemistry: nucleic acid structure and chemistry, denaturation of DNA, purification and detection of DNA, DNA electrophoresis, restriction endonucleases, DNA cloning, cloning vectors, genetic engineering, polymerase chain reaction, DNA sequencing, protein synthesis, the genetic code. Metabolic Biochemistry: bioenergetics and biochemical reaction types, anabolism, catabolism, glycolysis, fermentation, gluconeogenesis, the pentose phosphate pathway, the citric acid cycle, electron-transfer reaction, oxidative phosphorylation, photophosphorylation, photosynthesis."
$ws.Cells.Item(16,3).Value = "Basic biochemistry: Solvent properties of water, acids and bases, titration curves, buffer solutions. Aminoacids: three-dimensional structure, structure-property relationship, isoelectric point, electrophoresis, peptide bond, small peptides with physiological activity. Proteins: structure levels, structural irregularities, denaturation and renaturation, isolation, characterization and quantification. Enzymes: general concepts and mechanisms of action, cofactors and coenzymes, kinetics of michaelian enzymes, equilibrium and reaction rate, kinetic parameters and their applications, allosteric enzymes, regulation and inhibition of enzymes, general mechanisms of enzymatic reactions. Carbohydrates: classification of monosaccharides, cyclic structure and isomerism, chemical properties, disaccharides, homopolysaccharides and their structural and energy storage functions, heteropolysaccharides, glycoproteins and glycolipids. Lipids: fatty acids, triacylglycerides, phospholipids, sphingolipids and cholesterol. Biological membranes: fluid mosaic model, transport, selective permeability, passive process and active transport of biomolecules and/or ions. Molecular Biochemistry: nucleic acid structure and chemistry, denaturation of DNA, purification and detection of DNA, DNA electrophoresis, restriction endonucleases, DNA cloning, cloning vectors, genetic engineering, polymerase chain reaction, DNA sequencing, protein synthesis, the genetic code. Metabolic Biochemistry: bioenergetics and biochemical reaction types, anabolism, catabolism, glycolysis, fermentation, gluconeogenesis, the pentose phosphate pathway, the citric acid cycle, electron-transfer reaction, oxidative phosphorylation, photophosphorylation, photosynthesis."
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Cells.Item(17,1).Value = "Avaliação:"
$ws.Cells.Item(17,2).ClearContents()
$ws.Cells.Item(17,3).ClearContents()

# Row 18
$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "5111420 - Talita Martins Lacerda"
$ws.Cells.Item(18,3).Value = "5111420 - Talita Martins Lacerda"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(19,2).Value = "Notas N distribuído no semestre. A composição das `"N`" fica critério do docente."
$ws.Cells.Item(19,3).Value = "Notas N distribuído no semestre. A composição das `"N`" fica critério do docente."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(20,2).Value = "MF = média aritmética ou ponderada das notas das avaliações (a critério do docente)"
$ws.Cells.Item(20,3).Value = "MF = média aritmética ou ponderada das notas das avaliações (a critério do docente)"
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(21,2).Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maiordo que 5,0."
$ws.Cells.Item(21,3).Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maiordo que 5,0."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Cells.Item(22,1).Value = "Requisitos:"
$ws.Cells.Item(22,2).ClearContents()
$ws.Cells.Item(22,3).ClearContents()

# Row 23
$ws.Cells.Item(23,1).ClearContents()
$ws.Cells.Item(23,2).Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Cells.Item(23,3).Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("A1").Select()